# "Generate Report for Handback" — for each localized-language sheet
# (zh-cn, de-de), mark the two handed-off rows as handed back:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Target File" / "Latest Handback File" columns (E/F) get
#     populated with hyperlinks to the files that went out/came back
#   - "Latest Handback DateTime" (G) is stamped with the handback time

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$sheetsInfo = @(
    @{
        Name            = "zh-cn"
        MdUrl           = "https://github.com/OpenLocalizationTest/oltest/blob/2fd0b40863ecda8fed10c88a6c8e007ba23afa6e/e2e/1214c0cf-72f0-40f9-ad5a-bacb0e69c537.md"
        MdDisplay       = "1214c0cf-72f0-40f9-ad5a-bacb0e69c537.md"
        XlfUrl          = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/81d8065737627ac49bca568989ed737a0be011d5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/1214c0cf-72f0-40f9-ad5a-bacb0e69c537.9202164902571d7a033129dc8e57e3d8a2cd1b5d.zh-cn.xlf"
        XlfDisplay      = "1214c0cf-72f0-40f9-ad5a-bacb0e69c537.9202164902571d7a033129dc8e57e3d8a2cd1b5d.zh-cn.xlf"
        HandbackDateTime = "2016-03-07 02:40:02"
    },
    @{
        Name            = "de-de"
        MdUrl           = "https://github.com/OpenLocalizationTest/oltest/blob/2fd0b40863ecda8fed10c88a6c8e007ba23afa6e/e2e/1214c0cf-72f0-40f9-ad5a-bacb0e69c537.md"
        MdDisplay       = "1214c0cf-72f0-40f9-ad5a-bacb0e69c537.md"
        XlfUrl          = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/31a24e9af61abc576a200769325174c7b840bb9a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/1214c0cf-72f0-40f9-ad5a-bacb0e69c537.9202164902571d7a033129dc8e57e3d8a2cd1b5d.de-de.xlf"
        XlfDisplay      = "1214c0cf-72f0-40f9-ad5a-bacb0e69c537.9202164902571d7a033129dc8e57e3d8a2cd1b5d.de-de.xlf"
        HandbackDateTime = "2016-03-07 02:40:21"
    }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    foreach ($row in @(2, 3)) {
        # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
        $ws.Cells.Item($row, 2).Value = $newStatus

        # E = Latest Target File, F = Latest Handback File
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.Value = $info.MdDisplay
        $ws.Hyperlinks.Add($eCell, $info.MdUrl, "", "", $info.MdDisplay) | Out-Null

        $fCell = $ws.Cells.Item($row, 6)
        $fCell.Value = $info.XlfDisplay
        $ws.Hyperlinks.Add($fCell, $info.XlfUrl, "", "", $info.XlfDisplay) | Out-Null

        # G = Latest Handback DateTime
        $ws.Cells.Item($row, 7).Value = $info.HandbackDateTime
    }
}

Write-Host "Handback report generated"
